$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits right after the
#    "Referirse a" run (before " CPF2") in the first table.
#    The COM bridge doesn't expose Bookmarks.Delete/Item, but doing a
#    Find & Replace across the bookmark's position removes it as a
#    natural side effect of retyping that span of text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Referirse a CPF2", $true, $false, $false, $false, $false,
    $true, 1, $false, "Referirse a CPF2", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Fix the misused verb in the EF-3 step: "indica" -> "almacena",
#    which also relocates the "_GoBack" bookmark into the middle of
#    this sentence (right after "almacena ").
# ------------------------------------------------------------------
$oldText = "El sistema indica la cantidad de locs agregados, borrados y totales del programa."
$newText = "El sistema almacena  la cantidad de locs agregados, borrados y totales del programa."

$d.Content.Find.Execute(
    $oldText, $true, $false, $false, $false, $false,
    $true, 1, $false, $newText, 2) | Out-Null

# Locate the boundary right after "El sistema almacena " (and before
# " la cantidad ...") so the new bookmark lands between "almacena "
# and the rest of the sentence, matching the edited document.
$locate = $d.Content
$locate.Find.Execute("El sistema almacena ", $true) | Out-Null
$bmStart = $locate.End

$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
